$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.622.40'
$ws.Range('E2').Value = '  +2.24%  '
$ws.Range('D3').Value = '2.155.70'
$ws.Range('E3').Value = '  +2.73%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'226.86"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('E6').Value = '  +1.78%  '
$ws.Range('D7').Value = "'63.06"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.00%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = "'0.390"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('D10').Value = "'0.0844"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').Value = "'15.88"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('D13').Value = '2.475.85'
$ws.Range('E13').Value = '  +2.75%  '
$ws.Range('D14').Value = "'21.84"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '2.156.44'
$ws.Range('E17').Value = '  +2.66%  '
$ws.Range('D18').Value = '39.603.77'
$ws.Range('E18').Value = '  +2.19%  '
$ws.Range('D19').Value = "'71.55"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'6.06"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = "'229.47"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.14%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = "'2.40"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.72%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = "'2.36"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.41%  '
$ws.Range('D26').Value = "'172.06"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('D27').Value = "'9.54"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('E28').Value = '  +2.43%  '
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').Value = "'19.77"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.54%  '
$ws.Range('D31').Value = "'2.68"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.46%  '
$ws.Range('E32').Value = '  +1.50%  '
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').Value = "'4.69"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = "'0.0618"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('B36').Value = 'THORChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D36').Value = "'6.90"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.36%  '
$ws.Range('D37').Value = "'2.40"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('D38').Value = "'3.64"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.25%  '
$ws.Range('D39').Value = "'5.08"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +22.43%  '
$ws.Range('D41').Value = "'102.86"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').Value = "'17.62"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D44').Value = '1.516.21'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = "'1.20"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = "'0.0917"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').Value = "'7.76"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').Value = "'50.25"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.75%  '
$ws.Range('E51').Value = '  +1.09%  '
